# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) to the affected Leve rows across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 599.7692
$ws.Range("I92").Value = 548
$ws.Range("J92").Value = 996.6667
$ws.Range("K92").Value = 548
$ws.Range("L92").Value = 996.6667
$ws.Range("M92").Value = 700
$ws.Range("N92").Value = -3492.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1514.3077
$ws.Range("I2").Value = 757.86206
$ws.Range("J2").Value = 3708
$ws.Range("K2").Value = 757.86206
$ws.Range("L2").Value = 3708
$ws.Range("M2").Value = -644.86206
$ws.Range("N2").Value = -3934

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 293841.22
$ws.Range("I32").Value = 3954.6125
$ws.Range("J32").Value = 1514416.5
$ws.Range("K32").Value = 3954.6125
$ws.Range("L32").Value = 1514416.5
$ws.Range("M32").Value = -3667.6125
$ws.Range("N32").Value = -1514990.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2746.975
$ws.Range("I61").Value = 2749.4167
$ws.Range("J61").Value = 2725
$ws.Range("K61").Value = 2749.4167
$ws.Range("L61").Value = 2725
$ws.Range("M61").Value = -2537.4167
$ws.Range("N61").Value = -3149

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1373.9231
$ws.Range("I110").Value = 1082.2727
$ws.Range("J110").Value = 2978
$ws.Range("K110").Value = 1082.2727
$ws.Range("L110").Value = 2978
$ws.Range("M110").Value = 962.7273
$ws.Range("N110").Value = -7068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1514.3077
$ws.Range("I116").Value = 757.86206
$ws.Range("J116").Value = 3708
$ws.Range("K116").Value = 757.86206
$ws.Range("L116").Value = 3708
$ws.Range("M116").Value = 1536.13794
$ws.Range("N116").Value = -8296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2746.975
$ws.Range("I136").Value = 2749.4167
$ws.Range("J136").Value = 2725
$ws.Range("K136").Value = 8248.250100000001
$ws.Range("L136").Value = 8175
$ws.Range("M136").Value = -5698.250100000001
$ws.Range("N136").Value = -13275

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1514.3077
$ws.Range("I3").Value = 757.86206
$ws.Range("J3").Value = 3708
$ws.Range("K3").Value = 757.86206
$ws.Range("L3").Value = 3708
$ws.Range("M3").Value = -643.86206
$ws.Range("N3").Value = -3936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2255.4773
$ws.Range("I20").Value = 1633.4073
$ws.Range("J20").Value = 3243.4707
$ws.Range("K20").Value = 1633.4073
$ws.Range("L20").Value = 3243.4707
$ws.Range("M20").Value = -1386.4073
$ws.Range("N20").Value = -3737.4707

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 250911.58
$ws.Range("I107").Value = 788.5161000000001
$ws.Range("J107").Value = 1112446.5
$ws.Range("K107").Value = 788.5161000000001
$ws.Range("L107").Value = 1112446.5
$ws.Range("M107").Value = 1131.4839
$ws.Range("N107").Value = -1116286.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 50390
$ws.Range("J132").Value = 50390
$ws.Range("L132").Value = 50390
$ws.Range("N132").Value = -60510

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1828.4872
$ws.Range("I134").Value = 1188.2963
$ws.Range("J134").Value = 3268.9167
$ws.Range("K134").Value = 3564.8889
$ws.Range("L134").Value = 9806.750100000001
$ws.Range("M134").Value = -1029.8889
$ws.Range("N134").Value = -14876.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 56720
$ws.Range("J138").Value = 56720
$ws.Range("L138").Value = 56720
$ws.Range("N138").Value = -67000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 566.53845
$ws.Range("I16").Value = 586.7
$ws.Range("J16").Value = 499.33334
$ws.Range("K16").Value = 586.7
$ws.Range("L16").Value = 499.33334
$ws.Range("M16").Value = -299.7
$ws.Range("N16").Value = -1073.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 566.53845
$ws.Range("I113").Value = 586.7
$ws.Range("J113").Value = 499.33334
$ws.Range("K113").Value = 586.7
$ws.Range("L113").Value = 499.33334
$ws.Range("M113").Value = 1583.3
$ws.Range("N113").Value = -4839.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1713.0435
$ws.Range("I132").Value = 1205.2122
$ws.Range("J132").Value = 3002.1538
$ws.Range("K132").Value = 3615.6366
$ws.Range("L132").Value = 9006.4614
$ws.Range("M132").Value = -1085.6366
$ws.Range("N132").Value = -14066.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 894.5526
$ws.Range("I134").Value = 705.5357
$ws.Range("J134").Value = 1423.8
$ws.Range("K134").Value = 2116.6071
$ws.Range("L134").Value = 4271.4
$ws.Range("M134").Value = 418.3928999999998
$ws.Range("N134").Value = -9341.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 88500
$ws.Range("J140").Value = 88500
$ws.Range("L140").Value = 88500
$ws.Range("N140").Value = -98860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 4000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 12000
$ws.Range("M104").ClearContents() | Out-Null
$ws.Range("N104").Value = -17242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1149.409
$ws.Range("I114").Value = 1004.0714
$ws.Range("J114").Value = 1403.75
$ws.Range("K114").Value = 3012.2142
$ws.Range("L114").Value = 4211.25
$ws.Range("M114").Value = 241.7857999999997
$ws.Range("N114").Value = -10719.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15210124
$ws.Range("I70").Value = 24463056
$ws.Range("J70").Value = 8879.214
$ws.Range("K70").Value = 24463056
$ws.Range("L70").Value = 8879.214
$ws.Range("M70").Value = -24462786
$ws.Range("N70").Value = -9419.214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 15210124
$ws.Range("I73").Value = 24463056
$ws.Range("J73").Value = 8879.214
$ws.Range("K73").Value = 24463056
$ws.Range("L73").Value = 8879.214
$ws.Range("M73").Value = -24462120
$ws.Range("N73").Value = -10751.214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1963.3954
$ws.Range("I132").Value = 1581.3715
$ws.Range("J132").Value = 3634.75
$ws.Range("K132").Value = 4744.1145
$ws.Range("L132").Value = 10904.25
$ws.Range("M132").Value = -2214.1145
$ws.Range("N132").Value = -15964.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3000
$ws.Range("J20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3452

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2254.2727
$ws.Range("I40").Value = 2112.7856
$ws.Range("J40").Value = 2501.875
$ws.Range("K40").Value = 2112.7856
$ws.Range("L40").Value = 2501.875
$ws.Range("M40").Value = -1976.7856
$ws.Range("N40").Value = -2773.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 40499.375
$ws.Range("J47").Value = 40499.375
$ws.Range("L47").Value = 40499.375
$ws.Range("N47").Value = -41479.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 40499.375
$ws.Range("J52").Value = 40499.375
$ws.Range("L52").Value = 40499.375
$ws.Range("N52").Value = -40965.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1570.9387
$ws.Range("I93").Value = 1083.1875
$ws.Range("J93").Value = 2489.0588
$ws.Range("K93").Value = 1083.1875
$ws.Range("L93").Value = 2489.0588
$ws.Range("M93").Value = 164.8125
$ws.Range("N93").Value = -4985.0588

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2999.5293
$ws.Range("I122").Value = 2734.25
$ws.Range("J122").Value = 3636.2
$ws.Range("K122").Value = 8202.75
$ws.Range("L122").Value = 10908.6
$ws.Range("M122").Value = -5752.75
$ws.Range("N122").Value = -15808.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 25000514
$ws.Range("I113").Value = 457.63635
$ws.Range("J113").Value = 142857920
$ws.Range("K113").Value = 1372.90905
$ws.Range("L113").Value = 428573760
$ws.Range("M113").Value = 797.09095
$ws.Range("N113").Value = -428578100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 28847520
$ws.Range("I132").Value = 38462644
$ws.Range("J132").Value = 2152.7693
$ws.Range("K132").Value = 115387932
$ws.Range("L132").Value = 6458.3079
$ws.Range("M132").Value = -115385402
$ws.Range("N132").Value = -11518.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 698.95917
$ws.Range("I136").Value = 579.8889
$ws.Range("J136").Value = 1028.6923
$ws.Range("K136").Value = 1739.6667
$ws.Range("L136").Value = 3086.0769
$ws.Range("M136").Value = 810.3332999999998
$ws.Range("N136").Value = -8186.0769
